$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '56.024.31'
$ws.Range('E2').Value = '  -3.72%  '
$ws.Range('D3').Value = '2.362.28'
$ws.Range('E3').Value = '  -4.05%  '
$ws.Range('E4').Value = '  -0.05%  '
$ws.Range('D5').Value = "'502.04"
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -1.58%  '
$ws.Range('D6').Value = "'129.40"
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -3.47%  '
$ws.Range('E7').Value = '  -0.10%  '
$ws.Range('D8').Value = "'0.544"
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  -2.88%  '
$ws.Range('D9').Value = '2.364.06'
$ws.Range('E9').Value = '  -4.06%  '
$ws.Range('D10').Value = "'0.0985"
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +0.06%  '
$ws.Range('E11').Value = '  -0.06%  '
$ws.Range('D12').Value = "'4.81"
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  +3.17%  '
$ws.Range('D13').Value = "'0.325"
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  +0.10%  '
$ws.Range('D14').Value = '2.780.90'
$ws.Range('E14').Value = '  -4.21%  '
$ws.Range('D15').Value = '55.945.23'
$ws.Range('E15').Value = '  -3.65%  '
$ws.Range('D16').Value = "'21.42"
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  -2.82%  '
$ws.Range('E17').Value = '  -3.35%  '
$ws.Range('D18').Value = '2.362.10'
$ws.Range('E18').Value = '  -4.12%  '
$ws.Range('D19').Value = "'10.01"
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -3.57%  '
$ws.Range('B20').Value = 'Polkadot'
$ws.Range('C20').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range('D20').Value = "'4.01"
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -2.68%  '
$ws.Range('B21').Value = 'BitcoinCash'
$ws.Range('C21').Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range('D21').Value = "'307.57"
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -2.67%  '
$ws.Range('E22').Value = '  -2.29%  '
$ws.Range('D23').Value = "'1.00"
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +0.12%  '
$ws.Range('D24').Value = "'65.86"
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +0.41%  '
$ws.Range('E25').Value = '  -0.22%  '
$ws.Range('D26').Value = "'0.369"
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -3.18%  '
$ws.Range('E27').Value = '  -6.33%  '
$ws.Range('D28').Value = "'7.24"
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -5.02%  '
$ws.Range('D29').Value = "'172.77"
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +0.72%  '
$ws.Range('E30').Value = '  -3.82%  '
$ws.Range('E31').Value = '  -3.03%  '
$ws.Range('E33').Value = '  -6.04%  '
$ws.Range('E34').Value = '  -0.30%  '
$ws.Range('E35').Value = '  -5.84%  '
$ws.Range('D36').Value = "'17.62"
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -2.74%  '
$ws.Range('E37').Value = '  -5.99%  '
$ws.Range('D38').Value = "'3.72"
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -4.85%  '
$ws.Range('D39').Value = "'36.15"
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -1.93%  '
$ws.Range('D40').Value = "'0.801"
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -1.43%  '
$ws.Range('D41').Value = "'1.38"
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -6.20%  '
$ws.Range('E42').Value = '  -1.38%  '
$ws.Range('D43').Value = "'128.82"
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -5.98%  '
$ws.Range('D44').Value = "'4.67"
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -5.37%  '
$ws.Range('D45').Value = "'0.561"
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -2.73%  '
$ws.Range('D46').Value = "'0.0904"
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -1.60%  '
$ws.Range('D47').Value = "'238.16"
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -7.24%  '
$ws.Range('E48').Value = '  -2.77%  '
$ws.Range('E49').Value = '  -3.98%  '
$ws.Range('D50').Value = "'16.97"
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -2.38%  '
$ws.Range('E51').Value = '  -1.50%  '
